# Update "想去人数" (F column) figures for several events on sheets
# "展览" and "全部类型" to the newly scraped counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1834
$ws1.Range("F9").Value  = 2425
$ws1.Range("F13").Value = 1468
$ws1.Range("F14").Value = 515
$ws1.Range("F15").Value = 36
$ws1.Range("F17").Value = 224
$ws1.Range("F22").Value = 213
$ws1.Range("F24").Value = 119
$ws1.Range("F29").Value = 335
$ws1.Range("F30").Value = 186

# Sheet "全部类型": row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1834
$ws4.Range("F10").Value = 2425
$ws4.Range("F14").Value = 1468
$ws4.Range("F15").Value = 515
$ws4.Range("F16").Value = 36
$ws4.Range("F18").Value = 224
$ws4.Range("F23").Value = 213
$ws4.Range("F25").Value = 119
$ws4.Range("F30").Value = 335
$ws4.Range("F31").Value = 186
